{"js": "// Replacement pairs in document order: the date heading, then each\n// table cell's math expression (row-major, left-to-right, top-to-bottom).\nconst replacements = [\n  [\"2024-05-31 Friday\", \"2024-06-01 Saturday\"],\n  [\"68-66=\", \"57+24=\"],\n  [\"83-45=\", \"33+42=\"],\n  [\"69-15=\", \"92-49=\"],\n  [\"58+5=\", \"24-7=\"],\n  [\"46-5=\", \"55-29=\"],\n  [\"63-38=\", \"38+24=\"],\n  [\"33+62=\", \"19+71=\"],\n  [\"0+49=\", \"40+13=\"],\n  [\"78+0=\", \"33+43=\"],\n  [\"6+27=\", \"56+9=\"],\n  [\"59-5=\", \"43-14=\"],\n  [\"30+25=\", \"85-82=\"],\n  [\"6+41=\", \"59+29=\"],\n  [\"57-54=\", \"4+60=\"],\n  [\"78+7=\", \"91-4=\"],\n  [\"4+52=\", \"79+14=\"],\n  [\"25+55=\", \"11+37=\"],\n  [\"10+66=\", \"43+50=\"],\n  [\"90-8=\", \"2+38=\"],\n  [\"20+60=\", \"72-40=\"],\n  [\"70+5=\", \"1+20=\"],\n  [\"56-7=\", \"54-18=\"],\n  [\"10+9=\", \"51+31=\"],\n  [\"21+63=\", \"97-76=\"],\n  [\"44+42=\", \"92-3=\"],\n  [\"16+28=\", \"63-8=\"],\n  [\"36+46=\", \"43+44=\"],\n  [\"77-33=\", \"90-38=\"],\n  [\"16+18=\", \"22+37=\"],\n  [\"51+9=\", \"54-51=\"],\n  [\"72+15=\", \"74-61=\"],\n  [\"0+93=\", \"17+47=\"],\n  [\"63-11=\", \"38+3=\"],\n  [\"21-19=\", \"26+68=\"],\n  [\"38+56=\", \"41+57=\"],\n  [\"45-20=\", \"75-68=\"],\n  [\"26+24=\", \"5-3=\"],\n  [\"62-34=\", \"33+57=\"],\n  [\"26+69=\", \"9+63=\"],\n  [\"10+60=\", \"70-50=\"],\n  [\"63-47=\", \"93-22=\"],\n  [\"74-44=\", \"49+37=\"],\n  [\"72+16=\", \"19+30=\"],\n  [\"83-76=\", \"56-44=\"],\n  [\"86-83=\", \"14+68=\"],\n  [\"95-38=\", \"92-53=\"],\n  [\"55-6=\", \"45-17=\"],\n  [\"41-32=\", \"88-20=\"],\n  [\"54-43=\", \"49-17=\"],\n  [\"70-44=\", \"47+45=\"],\n  [\"8+4=\", \"57-49=\"],\n  [\"93-14=\", \"85-24=\"],\n  [\"27-26=\", \"44-25=\"],\n  [\"29+42=\", \"63-43=\"],\n  [\"13+35=\", \"47-39=\"],\n  [\"70-19=\", \"74+11=\"],\n  [\"36+19=\", \"76+19=\"],\n  [\"91-67=\", \"19+23=\"],\n  [\"44+30=\", \"97-56=\"],\n  [\"0+46=\", \"52+3=\"],\n  [\"14+15=\", \"61-40=\"],\n  [\"28-26=\", \"2+58=\"],\n  [\"73+22=\", \"69-11=\"],\n  [\"13+51=\", \"70-21=\"],\n  [\"31+48=\", \"1+32=\"],\n  [\"67-32=\", \"23+8=\"],\n  [\"11+41=\", \"51+4=\"],\n  [\"60-29=\", \"57-20=\"],\n  [\"92-52=\", \"13+31=\"],\n  [\"1+64=\", \"5+8=\"],\n  [\"0+69=\", \"81-47=\"],\n  [\"8+80=\", \"17+9=\"],\n  [\"27-6=\", \"61-35=\"],\n  [\"73-39=\", \"79-65=\"],\n  [\"48-14=\", \"63+12=\"],\n  [\"33-30=\", \"56-2=\"],\n  [\"60-9=\", \"72-31=\"],\n  [\"87-12=\", \"3+11=\"],\n  [\"50-46=\", \"23+48=\"],\n  [\"83-53=\", \"65+11=\"],\n  [\"75-72=\", \"85-49=\"],\n  [\"48-5=\", \"82-38=\"],\n  [\"48+11=\", \"52+38=\"],\n  [\"35+16=\", \"4+90=\"],\n  [\"64-13=\", \"45-1=\"],\n  [\"27-19=\", \"25+12=\"],\n  [\"53+15=\", \"52-51=\"],\n  [\"71-16=\", \"65+4=\"],\n  [\"8+18=\", \"48+23=\"],\n  [\"15+77=\", \"35+25=\"],\n  [\"54+5=\", \"14+36=\"],\n  [\"36-10=\", \"92-50=\"],\n  [\"8+46=\", \"84-62=\"],\n  [\"8+62=\", \"32+40=\"],\n  [\"37-13=\", \"83+14=\"],\n  [\"85+5=\", \"55+5=\"],\n  [\"99-12=\", \"7+45=\"],\n  [\"97-74=\", \"67+14=\"],\n  [\"87-45=\", \"26+42=\"],\n  [\"1+37=\", \"76-29=\"],\n];\n\nconst body = context.document.body;\n\n// Phase 1: swap every old string for a unique, collision-proof placeholder.\n// (Some new values are substrings of other old values -- e.g. \"11+37=\" contains\n// \"1+37=\" -- so a direct single-pass old->new replace could over-match text that\n// a previous iteration just inserted. Routing through placeholders first avoids that.)\nfor (let i = 0; i < replacements.length; i++) {\n  const oldText = replacements[i][0];\n  const placeholder = '\\u0001P' + i + '\\u0002';\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let k = 0; k < results.items.length; k++) {\n    results.items[k].insertText(placeholder, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// Phase 2: swap every placeholder for its final new value.\nfor (let i = 0; i < replacements.length; i++) {\n  const newText = replacements[i][1];\n  const placeholder = '\\u0001P' + i + '\\u0002';\n  const results = body.search(placeholder, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let k = 0; k < results.items.length; k++) {\n    results.items[k].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "# Replacement pairs in document order: the date heading, then each table\n# cell's math expression (row-major, left-to-right, top-to-bottom).\n$replacements = @(\n  ,@(\"2024-05-31 Friday\", \"2024-06-01 Saturday\")\n  ,@(\"68-66=\", \"57+24=\")\n  ,@(\"83-45=\", \"33+42=\")\n  ,@(\"69-15=\", \"92-49=\")\n  ,@(\"58+5=\", \"24-7=\")\n  ,@(\"46-5=\", \"55-29=\")\n  ,@(\"63-38=\", \"38+24=\")\n  ,@(\"33+62=\", \"19+71=\")\n  ,@(\"0+49=\", \"40+13=\")\n  ,@(\"78+0=\", \"33+43=\")\n  ,@(\"6+27=\", \"56+9=\")\n  ,@(\"59-5=\", \"43-14=\")\n  ,@(\"30+25=\", \"85-82=\")\n  ,@(\"6+41=\", \"59+29=\")\n  ,@(\"57-54=\", \"4+60=\")\n  ,@(\"78+7=\", \"91-4=\")\n  ,@(\"4+52=\", \"79+14=\")\n  ,@(\"25+55=\", \"11+37=\")\n  ,@(\"10+66=\", \"43+50=\")\n  ,@(\"90-8=\", \"2+38=\")\n  ,@(\"20+60=\", \"72-40=\")\n  ,@(\"70+5=\", \"1+20=\")\n  ,@(\"56-7=\", \"54-18=\")\n  ,@(\"10+9=\", \"51+31=\")\n  ,@(\"21+63=\", \"97-76=\")\n  ,@(\"44+42=\", \"92-3=\")\n  ,@(\"16+28=\", \"63-8=\")\n  ,@(\"36+46=\", \"43+44=\")\n  ,@(\"77-33=\", \"90-38=\")\n  ,@(\"16+18=\", \"22+37=\")\n  ,@(\"51+9=\", \"54-51=\")\n  ,@(\"72+15=\", \"74-61=\")\n  ,@(\"0+93=\", \"17+47=\")\n  ,@(\"63-11=\", \"38+3=\")\n  ,@(\"21-19=\", \"26+68=\")\n  ,@(\"38+56=\", \"41+57=\")\n  ,@(\"45-20=\", \"75-68=\")\n  ,@(\"26+24=\", \"5-3=\")\n  ,@(\"62-34=\", \"33+57=\")\n  ,@(\"26+69=\", \"9+63=\")\n  ,@(\"10+60=\", \"70-50=\")\n  ,@(\"63-47=\", \"93-22=\")\n  ,@(\"74-44=\", \"49+37=\")\n  ,@(\"72+16=\", \"19+30=\")\n  ,@(\"83-76=\", \"56-44=\")\n  ,@(\"86-83=\", \"14+68=\")\n  ,@(\"95-38=\", \"92-53=\")\n  ,@(\"55-6=\", \"45-17=\")\n  ,@(\"41-32=\", \"88-20=\")\n  ,@(\"54-43=\", \"49-17=\")\n  ,@(\"70-44=\", \"47+45=\")\n  ,@(\"8+4=\", \"57-49=\")\n  ,@(\"93-14=\", \"85-24=\")\n  ,@(\"27-26=\", \"44-25=\")\n  ,@(\"29+42=\", \"63-43=\")\n  ,@(\"13+35=\", \"47-39=\")\n  ,@(\"70-19=\", \"74+11=\")\n  ,@(\"36+19=\", \"76+19=\")\n  ,@(\"91-67=\", \"19+23=\")\n  ,@(\"44+30=\", \"97-56=\")\n  ,@(\"0+46=\", \"52+3=\")\n  ,@(\"14+15=\", \"61-40=\")\n  ,@(\"28-26=\", \"2+58=\")\n  ,@(\"73+22=\", \"69-11=\")\n  ,@(\"13+51=\", \"70-21=\")\n  ,@(\"31+48=\", \"1+32=\")\n  ,@(\"67-32=\", \"23+8=\")\n  ,@(\"11+41=\", \"51+4=\")\n  ,@(\"60-29=\", \"57-20=\")\n  ,@(\"92-52=\", \"13+31=\")\n  ,@(\"1+64=\", \"5+8=\")\n  ,@(\"0+69=\", \"81-47=\")\n  ,@(\"8+80=\", \"17+9=\")\n  ,@(\"27-6=\", \"61-35=\")\n  ,@(\"73-39=\", \"79-65=\")\n  ,@(\"48-14=\", \"63+12=\")\n  ,@(\"33-30=\", \"56-2=\")\n  ,@(\"60-9=\", \"72-31=\")\n  ,@(\"87-12=\", \"3+11=\")\n  ,@(\"50-46=\", \"23+48=\")\n  ,@(\"83-53=\", \"65+11=\")\n  ,@(\"75-72=\", \"85-49=\")\n  ,@(\"48-5=\", \"82-38=\")\n  ,@(\"48+11=\", \"52+38=\")\n  ,@(\"35+16=\", \"4+90=\")\n  ,@(\"64-13=\", \"45-1=\")\n  ,@(\"27-19=\", \"25+12=\")\n  ,@(\"53+15=\", \"52-51=\")\n  ,@(\"71-16=\", \"65+4=\")\n  ,@(\"8+18=\", \"48+23=\")\n  ,@(\"15+77=\", \"35+25=\")\n  ,@(\"54+5=\", \"14+36=\")\n  ,@(\"36-10=\", \"92-50=\")\n  ,@(\"8+46=\", \"84-62=\")\n  ,@(\"8+62=\", \"32+40=\")\n  ,@(\"37-13=\", \"83+14=\")\n  ,@(\"85+5=\", \"55+5=\")\n  ,@(\"99-12=\", \"7+45=\")\n  ,@(\"97-74=\", \"67+14=\")\n  ,@(\"87-45=\", \"26+42=\")\n  ,@(\"1+37=\", \"76-29=\")\n)\n\n$d = $word.ActiveDocument\n\n# wdFindContinue / wdReplaceAll constants used with Find.Execute below.\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n# Phase 1: swap every old string for a unique, collision-proof placeholder.\n# (Some new values are substrings of other old values -- e.g. \"11+37=\" contains\n# \"1+37=\" -- so a direct single-pass old->new replace could over-match text that\n# a previous iteration just inserted. Routing through placeholders first avoids that.)\nfor ($i = 0; $i -lt $replacements.Count; $i++) {\n  $oldText = $replacements[$i][0]\n  $placeholder = [char]0x0001 + \"P$i\" + [char]0x0002\n  $rng = $d.Content\n  $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $placeholder, $wdReplaceAll) | Out-Null\n}\n\n# Phase 2: swap every placeholder for its final new value.\nfor ($i = 0; $i -lt $replacements.Count; $i++) {\n  $newText = $replacements[$i][1]\n  $placeholder = [char]0x0001 + \"P$i\" + [char]0x0002\n  $rng = $d.Content\n  $rng.Find.Execute($placeholder, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll) | Out-Null\n}"}
